$wb = $excel.ActiveWorkbook

# Replace the status text "Ready for handoff" -> "In Translation" everywhere it occurs.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C4").Value = "In Translation"

# Shrink the Status columns' width to match the now-shorter text
# (target stored width 13.4101845877511; the host quantizes ColumnWidth to
# 1/6-character steps, so 12.5 is the closest input that lands on it).
$overview.Range("E:F").ColumnWidth = 12.5
$zhcn.Range("C:C").ColumnWidth = 12.5
$dede.Range("C:C").ColumnWidth = 12.5
